$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume cells being updated so Excel keeps them as literal text
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "E47", "E48", "E49", "E50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "328.61"
$ws.Range("E2").Value = "0.36%"
$ws.Range("D3").Value = "44.05"
$ws.Range("E3").Value = "0.15%"
$ws.Range("D4").Value = "5.578"
$ws.Range("E4").Value = "1.82%"
$ws.Range("D5").Value = "0.08061"
$ws.Range("E5").Value = "-0.11%"
$ws.Range("D6").Value = "2.002"
$ws.Range("E6").Value = "6.47%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "0.9524"
$ws.Range("E7").Value = "1.79%"
$ws.Range("D8").Value = "2.560"
$ws.Range("E8").Value = "-6.01%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1159"
$ws.Range("E9").Value = "-0.87%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1856"
$ws.Range("E10").Value = "-1.81%"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").Value = "11.84"
$ws.Range("E11").Value = "38.39%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09870"
$ws.Range("E12").Value = "3.18%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.04747"
$ws.Range("E13").Value = "14.04%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.1068"
$ws.Range("E14").Value = "0.23%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001285"
$ws.Range("E15").Value = "1.39%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "0.04238"
$ws.Range("E16").Value = "-2.39%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.005879"
$ws.Range("E17").Value = "-2.03%"
$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D18").Value = "0.004321"
$ws.Range("E18").Value = "-1.22%"
$ws.Range("D19").Value = "3.370"
$ws.Range("E19").Value = "-5.66%"
$ws.Range("B20").Value = "GateToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D20").Value = "4.330"
$ws.Range("E20").Value = "1.14%"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "0.3475"
$ws.Range("E21").Value = "-0.29%"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "0.1410"
$ws.Range("E22").Value = "3.17%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.2509"
$ws.Range("E23").Value = "-3.30%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "0.001253"
$ws.Range("E24").Value = "1.48%"
$ws.Range("D25").Value = "0.0001193"
$ws.Range("E25").Value = "-3.16%"
$ws.Range("E26").Value = "-0.59%"
$ws.Range("E38").Value = "-0.47%"
$ws.Range("D39").Value = "0.05523"
$ws.Range("E39").Value = "1.48%"
$ws.Range("D40").Value = "0.007569"
$ws.Range("E40").Value = "-1.06%"
$ws.Range("D41").Value = "0.1404"
$ws.Range("E41").Value = "0.96%"
$ws.Range("D42").Value = "0.008087"
$ws.Range("E42").Value = "-29.36%"
$ws.Range("D43").Value = "0.002020"
$ws.Range("E43").Value = "-4.34%"
$ws.Range("D44").Value = "0.008382"
$ws.Range("E44").Value = "-12.89%"
$ws.Range("D45").Value = "0.00007098"
$ws.Range("E45").Value = "2.54%"
$ws.Range("E46").Value = "-0.19%"
$ws.Range("E47").Value = "1.07%"
$ws.Range("E48").Value = "36.26%"
$ws.Range("E49").Value = "-0.19%"
$ws.Range("E50").Value = "-0.19%"
